# Streak and points.pptx — "Added businessmodels, pitch tips"
#
# 1. Remove the now-redundant "streaks system" heading textbox.
# 2. Nudge the "Levelcurve durch Funktion:" caption and its chart picture
#    to their new positions.
# 3. Add two small bold "Vers.1:" / "Vers.2:" version labels above the
#    two level-curve comparison charts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Delete the "streaks system" textbox (id 5, "Textfeld 4")
$s.Shapes.Item("Textfeld 4").Delete()

# 2) Reposition the "Levelcurve durch Funktion:" label (id 24, "Textfeld 23")
$lvlLabel = $s.Shapes.Item("Textfeld 23")
$lvlLabel.Left = 385.3279577559055
$lvlLabel.Top = 56.20181102362205

# 2) Reposition the level-curve chart picture (id 28, "Grafik 27")
$lvlChart = $s.Shapes.Item("Grafik 27")
$lvlChart.Left = 385.6633858267717
$lvlChart.Top = 84.94299212598425

# 3) Add "Vers.1:" / "Vers.2:" version labels above the two charts.
# Duplicate the neighbouring "Levelcurve durch Funktion:" caption so the
# new textboxes inherit the same no-fill / autofit-wrap="none" styling,
# then reposition, resize, retext and bold them.
$v1 = $lvlLabel.Duplicate().Item(1)
$v1.Name = "Textfeld 1"
$v1.Left = 460.72716535433074
$v1.Top = 26.65771653543307
$v1.Width = 65.89228446456693
$v1.Height = 29.081259842519685
$v1.TextFrame.TextRange.Text = "Vers.1:"
$v1.TextFrame.TextRange.Font.Bold = $true

$v2 = $lvlLabel.Duplicate().Item(1)
$v2.Name = "Textfeld 2"
$v2.Left = 782.9327699055118
$v2.Top = 26.095434070866144
$v2.Width = 65.89228446456693
$v2.Height = 29.081259842519685
$v2.TextFrame.TextRange.Text = "Vers.2:"
$v2.TextFrame.TextRange.Font.Bold = $true

# 4) File this (only) slide into a new "Leveling und Exp curve" section,
# matching the section bookkeeping PowerPoint added alongside the edit.
if ($p.SectionProperties.Count -eq 0) {
    [void]$p.SectionProperties.AddSection(1, "Leveling und Exp curve")
}
